$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 3.5
$ws.Range("J4").Value = 3.45
$ws.Range("L4").Value = 2.67
$ws.Range("R4").Value = 2
$ws.Range("W4").Value = 11.5
$ws.Range("Z4").Value = 37
$ws.Range("AA4").Value = 23
$ws.Range("AB4").Value = 27
$ws.Range("AD4").Value = 7
$ws.Range("AH4").Value = 9.25
$ws.Range("AI4").Value = 11.5
$ws.Range("AM4").Value = 23
$ws.Range("AO4").Value = 15.5
$ws.Range("AP4").Value = 20
$ws.Range("AQ4").Value = 70
$ws.Range("AR4").Value = 90
$ws.Range("AS4").Value = 200
$ws.Range("AX4").Value = 10.75
$ws.Range("AY4").Value = 17
$ws.Range("AZ4").Value = 40
$ws.Range("BA4").Value = 65
